# Weekly update: insert a new week's worth of data (2 rows) for
# "Hortaliza, Femacal de La Calera - Pimiento" at the top of the data
# block that starts at row 987, pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 987..end down by two rows, inserting two blank rows.
$ws.Rows("987:988").Insert()

# --- New row 987 ---------------------------------------------------
$ws.Range("A987").Value = 3
$ws.Range("B987").Value = "Femacal de La Calera"
$ws.Range("C987").Value = "Coquimbo"
$ws.Range("D987").Value = 44516
$ws.Range("E987").Value = 5
$ws.Range("F987").Value = 100112002
$ws.Range("G987").Value = "Pimiento"
$ws.Range("H987").Value = "Zafiro rojo"
$ws.Range("I987").Value = "Primera"
$ws.Range("J987").Value = 105
$ws.Range("K987").Value = 42000
$ws.Range("L987").Value = 43000
$ws.Range("M987").Value = 42476
$ws.Range("N987").Value = "`$/caja 15 kilos"
$ws.Range("O987").Value = "Región de Arica y Parinacota"
$ws.Range("P987").Value = 2832
$ws.Range("Q987").Value = 15
$ws.Range("R987").Value = "Hortaliza"

# --- New row 988 ---------------------------------------------------
$ws.Range("A988").Value = 3
$ws.Range("B988").Value = "Femacal de La Calera"
$ws.Range("C988").Value = "Coquimbo"
$ws.Range("D988").Value = 44516
$ws.Range("E988").Value = 5
$ws.Range("F988").Value = 100112002
$ws.Range("G988").Value = "Pimiento"
$ws.Range("H988").Value = "Zafiro verde"
$ws.Range("I988").Value = "Primera"
$ws.Range("J988").Value = 35
$ws.Range("K988").Value = 30000
$ws.Range("L988").Value = 30000
$ws.Range("M988").Value = 30000
$ws.Range("N988").Value = "`$/caja 18 kilos"
$ws.Range("O988").Value = "Provincia de Quillota"
$ws.Range("P988").Value = 1667
$ws.Range("Q988").Value = 18
$ws.Range("R988").Value = "Hortaliza"
